$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.705.41"
$ws.Range("E2").Value = "  -0.86%  "

$ws.Range("D3").Value = "2.029.88"
$ws.Range("E3").Value = "  -0.95%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.606"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.28%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.92"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.81%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.373"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.57%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0826"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.04%  "

$ws.Range("E11").Value = "  +0.11%  "

$ws.Range("D12").Value = "2.335.48"
$ws.Range("E12").Value = "  -0.72%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.45"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.83%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.90%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.768"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.19%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.73%  "

$ws.Range("D17").Value = "2.028.30"
$ws.Range("E17").Value = "  -0.78%  "

$ws.Range("D18").Value = "37.681.88"
$ws.Range("E18").Value = "  -0.83%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.61%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.02%  "

$ws.Range("D21").Value = "0.0₃0817"
$ws.Range("E21").Value = "  -1.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.78"
$ws.Range("D22").Style = "Normal"

$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.38%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.61%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.75%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.31%  "

$ws.Range("E28").Value = "  -1.53%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.68"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.34%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.75%  "

$ws.Range("E31").Value = "  -0.26%  "

$ws.Range("E32").Value = "  +7.23%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.43%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0601"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.43%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.63%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.50"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.26%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.32"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.97%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.38"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.36%  "

$ws.Range("E39").Value = "  -0.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.00%  "

$ws.Range("D41").Value = "1.523.74"
$ws.Range("E41").Value = "  -0.86%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.77%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0214"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.27%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0904"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.10%  "

$ws.Range("E47").Value = "  -1.33%  "

$ws.Range("E48").Value = "  -0.83%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.83%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.13%  "

$ws.Range("D51").Value = "2.225.25"
$ws.Range("E51").Value = "  -0.73%  "
